$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ТО" (column C) figures for each sales rep row in Table1.
$ws.Range("C4").Value = 110
$ws.Range("C5").Value = 100
$ws.Range("C6").Value = 100
$ws.Range("C7").Value = 100
$ws.Range("C8").Value = 100
$ws.Range("C9").Value = 1001
$ws.Range("C10").Value = 1010
$ws.Range("C11").Value = 1

# Move the active selection to match the author's final cursor position.
$ws.Range("H15").Select()
